$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (these are numeric-looking strings kept as inlineStr text in the
    # source workbook), then the quotePrefix style that leaves behind on
    # the cell is cleared so the cell keeps its original, style-less look.
    $ws.Range($cell).Value2 = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

Set-TextCell "D2" ("25.868.83")
Set-TextCell "E2" ("  -1.45%  ")
Set-TextCell "D3" ("1.638.15")
Set-TextCell "D4" ("1.005")
Set-TextCell "E4" ("  -0.31%  ")
Set-TextCell "D5" ("215.28")
Set-TextCell "E5" ("  -0.67%  ")
Set-TextCell "D6" ("0.5025")
Set-TextCell "E6" ("  -2.39%  ")
Set-TextCell "E7" ("  -0.56%  ")
Set-TextCell "D8" ("0.2568")
Set-TextCell "E8" ("  -1.37%  ")
Set-TextCell "D9" ("0.06378")
Set-TextCell "E9" ("  -1.42%  ")
Set-TextCell "D10" ("19.66")
Set-TextCell "E10" ("  -1.59%  ")
Set-TextCell "D11" ("0.07745")
Set-TextCell "E11" ("  -1.26%  ")
Set-TextCell "B12" ("Polkadot")
Set-TextCell "C12" ("https://coinranking.com/coin/25W7FG7om+polkadot-dot")
Set-TextCell "D12" ("4.254")
Set-TextCell "E12" ("  -1.40%  ")
Set-TextCell "B13" ("WrappedEther")
Set-TextCell "C13" ("https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth")
Set-TextCell "D13" ("1.632.40")
Set-TextCell "E13" ("  -1.59%  ")
Set-TextCell "D14" ("1.863.18")
Set-TextCell "E14" ("  -1.23%  ")
Set-TextCell "D15" ("0.5457")
Set-TextCell "E15" ("  -1.61%  ")
Set-TextCell "D16" ("0.0" + [string]([char]0x2085) + "7894")
Set-TextCell "E16" ("  -1.98%  ")
Set-TextCell "D17" ("63.99")
Set-TextCell "E17" ("  -0.59%  ")
Set-TextCell "D18" ("25.868.64")
Set-TextCell "E18" ("  -1.49%  ")
Set-TextCell "E19" ("  -0.39%  ")
Set-TextCell "D20" ("202.40")
Set-TextCell "E20" ("  -4.04%  ")
Set-TextCell "D21" ("4.388")
Set-TextCell "E21" ("  -0.76%  ")
Set-TextCell "D22" ("9.891")
Set-TextCell "E22" ("  -2.20%  ")
Set-TextCell "D23" ("5.965")
Set-TextCell "D24" ("1.004")
Set-TextCell "E24" ("  -0.48%  ")
Set-TextCell "D25" ("1.908")
Set-TextCell "E25" ("  +8.17%  ")
Set-TextCell "D26" ("140.63")
Set-TextCell "E26" ("  -2.88%  ")
Set-TextCell "D27" ("0.1132")
Set-TextCell "E27" ("  -3.79%  ")
Set-TextCell "E28" ("  -1.36%  ")
Set-TextCell "D29" ("6.749")
Set-TextCell "E29" ("  -3.72%  ")
Set-TextCell "D30" ("1.241")
Set-TextCell "E30" ("  -0.50%  ")
Set-TextCell "D31" ("0.04970")
Set-TextCell "E31" ("  -2.88%  ")
Set-TextCell "D32" ("3.260")
Set-TextCell "E32" ("  -3.21%  ")
Set-TextCell "D33" ("3.191")
Set-TextCell "E34" ("  -1.44%  ")
Set-TextCell "D35" ("2.368")
Set-TextCell "E35" ("  +0.48%  ")
Set-TextCell "B36" ("MXToken")
Set-TextCell "C36" ("https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx")
Set-TextCell "D36" ("2.631")
Set-TextCell "E36" ("  -3.98%  ")
Set-TextCell "B37" ("ARBITRUM")
Set-TextCell "C37" ("https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb")
Set-TextCell "D37" ("0.8920")
Set-TextCell "E37" ("  -3.84%  ")
Set-TextCell "B38" ("ImmutableX")
Set-TextCell "C38" ("https://coinranking.com/coin/Z96jIvLU7+immutablex-imx")
Set-TextCell "D38" ("0.5607")
Set-TextCell "E38" ("  -2.39%  ")
Set-TextCell "B39" ("Maker")
Set-TextCell "C39" ("https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr")
Set-TextCell "D39" ("1.146.31")
Set-TextCell "E39" ("  -2.40%  ")
Set-TextCell "D40" ("0.01564")
Set-TextCell "E40" ("  -1.76%  ")
Set-TextCell "D41" ("1.002")
Set-TextCell "E41" ("  -0.66%  ")
Set-TextCell "D42" ("5.670")
Set-TextCell "E42" ("  -0.93%  ")
Set-TextCell "D43" ("99.62")
Set-TextCell "E43" ("  -0.90%  ")
Set-TextCell "D44" ("0.8062")
Set-TextCell "E44" ("  -2.39%  ")
Set-TextCell "D45" ("1.775.17")
Set-TextCell "E45" ("  -1.24%  ")
Set-TextCell "D46" ("0.0" + [string]([char]0x2088) + "117")
Set-TextCell "E46" ("  +1.16%  ")
Set-TextCell "D47" ("0.4543")
Set-TextCell "E47" ("  -0.41%  ")
Set-TextCell "D48" ("1.003")
Set-TextCell "E48" ("  -0.53%  ")
Set-TextCell "D49" ("54.88")
Set-TextCell "E49" ("  -1.19%  ")
Set-TextCell "D50" ("0.05059")
Set-TextCell "E50" ("  -0.47%  ")
Set-TextCell "E51" ("  -0.68%  ")
